$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (contaminants): zero out the comparison counts/percentages.
# B2/C2 historically hold numbers stored as text; prefix with an apostrophe
# so Excel keeps them as text instead of coercing to numeric, then reset the
# cell style back to Normal so no stray number-format style lingers.
$ws.Range("B2").Value = "'0"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "'0"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0

# Row 4 (flow_base_flow): update EDT_greater_than_RTT count and percentage
$ws.Range("B4").Value = "'114"
$ws.Range("B4").Style = "Normal"
$ws.Range("D4").Value = 0.9
